# "Generate Report for Handoff"
# The localization status moved from "In Translation" to "Ready for handoff";
# refresh the per-language status cells, the Xliff/handoff generation
# timestamps, and let the Status columns resize to fit the new, longer text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 08:56:10"

# --- zh-cn handoff detail sheet -------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 08:55:59"

# --- de-de handoff detail sheet -------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 08:56:10"

# --- Resize the Status columns to fit "Ready for handoff" -----------
$wsOverview.Range("E:F").ColumnWidth = 16.3333333333333
$wsZhCn.Range("C:C").ColumnWidth = 16.3333333333333
$wsDeDe.Range("C:C").ColumnWidth = 16.3333333333333
